# Apply updates to the "想去人数" (F) and "最低票价" (G) columns
# on the "展览" and "全部类型" worksheets, per the scraped data refresh.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 8605
$ws1.Range("G3").Value = 40
$ws1.Range("F5").Value = 87
$ws1.Range("F6").Value = 1376
$ws1.Range("F7").Value = 135
$ws1.Range("F10").Value = 9376
$ws1.Range("F11").Value = 152
$ws1.Range("F12").Value = 97
$ws1.Range("F15").Value = 356
$ws1.Range("F16").Value = 6355
$ws1.Range("F17").Value = 1063
$ws1.Range("F18").Value = 85
$ws1.Range("F20").Value = 133

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 8605
$ws4.Range("G3").Value = 40
$ws4.Range("F5").Value = 87
$ws4.Range("F6").Value = 1376
$ws4.Range("F7").Value = 135
$ws4.Range("F12").Value = 9376
$ws4.Range("F13").Value = 152
$ws4.Range("F14").Value = 97
$ws4.Range("F17").Value = 356
$ws4.Range("F18").Value = 6355
$ws4.Range("F19").Value = 1063
$ws4.Range("F20").Value = 85
$ws4.Range("F22").Value = 133
